# Bot data is now filled properly with storage folder and its size.
# Word's last-edit ("_GoBack") bookmark needs to move from the empty
# paragraph near the end of the document to wrap the newly-typed storage
# size figure "20 971 520 " (between the opening "(" and the word
# "bytes") in the "Maximum size of the stored data..." paragraph.
#
# Re-adding a bookmark under a name that already exists simply relocates
# it (exactly like real Word does), removing the old bookmarkStart/End
# pair and re-numbering every bookmark id in document order - which is
# precisely the behaviour captured by the target diff.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("20 971 520 ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the storage-size text to bookmark"
}

$d.Bookmarks.Add("_GoBack", $rng)
